$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-ambiguous value updates ---
$ws.Range('D2').Value = '56.736.26'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '2.332.63'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = '2.337.41'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E11').Value = '  +4.84%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').Value = '2.748.51'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '56.735.83'
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '2.332.76'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E20').Value = '  +1.84%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('E25').Value = '  +4.23%  '
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('E28').Value = '  +9.46%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').Value = '0.0₃0740'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('E37').Value = '  -1.45%  '
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('E44').Value = '  +4.44%  '
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  +0.90%  '
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('E50').Value = '  +6.18%  '
$ws.Range('E51').Value = '  +3.57%  '

# --- Numeric-looking price strings that must stay literal text ---
# Force text format so Excel doesn't coerce these to numbers (which would
# drop significant trailing zeros / change representation), then restore
# the default cell style so no stray formatting is introduced.
$numericTextCells = @('D5', 'D6', 'D10', 'D11', 'D12', 'D14', 'D20', 'D22', 'D24', 'D26', 'D27', 'D28', 'D29', 'D33', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D50', 'D51')
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = '@'
}
$ws.Range('D5').Value = '515.50'
$ws.Range('D6').Value = '134.81'
$ws.Range('D10').Value = '0.103'
$ws.Range('D11').Value = '5.38'
$ws.Range('D12').Value = '0.152'
$ws.Range('D14').Value = '23.91'
$ws.Range('D20').Value = '326.65'
$ws.Range('D22').Value = '6.62'
$ws.Range('D24').Value = '60.85'
$ws.Range('D26').Value = '0.999'
$ws.Range('D27').Value = '7.98'
$ws.Range('D28').Value = '1.29'
$ws.Range('D29').Value = '169.85'
$ws.Range('D33').Value = '18.50'
$ws.Range('D40').Value = '38.36'
$ws.Range('D41').Value = '144.00'
$ws.Range('D42').Value = '0.380'
$ws.Range('D43').Value = '3.62'
$ws.Range('D44').Value = '277.26'
$ws.Range('D45').Value = '5.17'
$ws.Range('D50').Value = '17.99'
$ws.Range('D51').Value = '17.49'
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = 'Normal'
}

